$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 536, shifting the existing rows 536:563 down to 537:564
$ws.Rows(536).Insert()

# Populate the new row with the weekly price-report entry
$ws.Range("A536").Value = 8
$ws.Range("B536").Value = "Terminal La Palmera de La Serena"
$ws.Range("C536").Value = "Coquimbo"
$ws.Range("D536").Value = 44939
$ws.Range("E536").Value = 4
$ws.Range("F536").Value = 100114001
$ws.Range("G536").Value = "Papa"
$ws.Range("H536").Value = "Cardinal"
$ws.Range("I536").Value = "1a (cosecha)"
$ws.Range("J536").Value = 2000
$ws.Range("K536").Value = 11800
$ws.Range("L536").Value = 12000
$ws.Range("M536").Value = 11900
$ws.Range("N536").Value = "$/saco 25 kilos"
$ws.Range("O536").Value = "Provincia del Elquí"
$ws.Range("P536").Value = 476
$ws.Range("Q536").Value = 25
$ws.Range("R536").Value = "Hortaliza"
